$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the MVF placeholder requirement rows with the real "extended features"
# requirement text (rows 9-11, column A).
$ws.Range("A9").Value = "Our system should source tournament and game data in real-time from APIs for tournaments where an API is available."
$ws.Range("A10").Value = "The system should have the ability to 'crowd-source' result information. i.e. Users can enter the score for a game and once a threshold of submitted results have been submitted, the results for the game will be updated and players scored on their picks"
$ws.Range("A11").Value = "Our system should have a suitable schedule and process for automatic backups of the database. This can be integrated as part of the Heroku platform we plan to deploy on."

# Remove the now-unneeded extended-feature placeholder rows (old rows 12-20).
# This shifts the trailing SUM row up to row 12 and auto-adjusts its formula
# from SUM(H4:H20) to SUM(H4:H11).
$ws.Rows("12:20").Delete()

# Narrow column A now that the remaining text is shorter.
$ws.Columns("A").ColumnWidth = 34.6

# Restore the previously-selected cell (now outside the shrunk data range).
[void]$ws.Range("F19").Select()
